# Update - GeometryShader - Specular Mapping
# Applies the cell-content edits reflected in the commit's xml diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Milestone marks (Roman numeral in column E, "X" completion mark in column F) ---

# Row 22: mark feature complete (F22 = "X"); E22 already "III"
$ws.Range("F22").Value = "X"

# Row 34: feature achieved on Milestone III
$ws.Range("E34").Value = "III"
$ws.Range("F34").Value = "X"

# Row 35: feature achieved on Milestone III
$ws.Range("E35").Value = "III"
$ws.Range("F35").Value = "X"

# Row 39: mark feature complete (F39 = "X"); E39 already "III"
$ws.Range("F39").Value = "X"

# Row 67: mark feature complete (F67 = "X"); E67 already "III"
$ws.Range("F67").Value = "X"

# Row 68: mark feature complete (F68 = "X"); E68 already "III"
$ws.Range("F68").Value = "X"

# Row 90: Effective Use of GIT - also mark Milestone III complete
$ws.Range("E90").Value = "X"

# Row 97: add project source citation URL
$ws.Range("A97").Value = "https://github.com/Microsoft/HoloLensCompanionKit/tree/master/RemotingHostSample/RemotingHostSampleShared"

# --- View state: update the active selection ---
$ws.Range("E63").Select()

$wb.Save()
